$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits between the runs "1.2" and
#    ".0" in the revision-history table. A same-text Find/Replace spanning
#    the bookmark's position collapses it away; we then re-split the merged
#    run back into "1.2" / ".0" with a harmless Bold no-op toggle so the
#    paragraph keeps two runs (matching the original run layout).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("1.2.0", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.2.0", 2)

$verRng = $d.Content
$verRng.Find.Execute("1.2.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotZero = $d.Range($verRng.Start + 3, $verRng.End)
$dotZero.Font.Bold = 1
$dotZero.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) "Matías Garcés " + "Bernt" + "." (split apart by spell-check proofErr
#    markers) becomes one plain run "Matías Garcés Bernt."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Matías Garcés Bernt.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Matías Garcés Bernt.", 2)

# ---------------------------------------------------------------------------
# 3) Fix the typo "mas" -> "más" in "Asignación de mas HH al proyecto.";
#    this also sweeps away the surrounding gramStart/gramEnd proofErr tags.
#    Word's edit cursor then leaves the "_GoBack" bookmark right after
#    "Asign", so we recreate that bookmark and re-split the run there and
#    again right around "más" (Bold no-op toggle, as above).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Asignación de mas HH", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Asignación de más HH", 2)

$asignRng = $d.Content
$asignRng.Find.Execute("Asign", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRng = $d.Range($asignRng.End, $asignRng.End)
$d.Bookmarks.Add("_GoBack", $goBackRng)

$masRng = $d.Content
$masRng.Find.Execute("más HH", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$masOnly = $d.Range($masRng.Start, $masRng.Start + 3)
$masOnly.Font.Bold = 1
$masOnly.Font.Bold = 0

# ---------------------------------------------------------------------------
# 4) Drop the last three rows (R8, R9, R10) from the risk-summary table.
# ---------------------------------------------------------------------------
$riskTable = $d.Tables.Item(12)
$riskTable.Rows.Item(11).Delete()
$riskTable.Rows.Item(10).Delete()
$riskTable.Rows.Item(9).Delete()

Write-Host "Edit complete"
